$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title
$ws.Range("B2").Value = "Bikerr Defect Tracking Report 1.3"

# Copy the formatting from an existing fully-formatted data row (row 13) down
# onto the new data rows (14-18) so the new rows look like proper table rows
# instead of the "blank" placeholder rows they used to be.
$ws.Range("C13:I13").Copy()
$ws.Range("C14:I18").PasteSpecial(-4122)
$ws.Range("J13").Copy()
$ws.Range("J14:J18").PasteSpecial(-4122)

# Row 14
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = "whole application"
$ws.Range("D14").Value = "Functional Error"
$ws.Range("E14").Value = "Major"
$ws.Range("F14").Value = "SE"
$ws.Range("G14").Value = "Upon token expiration, user is not logged out"
$ws.Range("H14").Value = "Sacha Elkaim"
$ws.Range("I14").Value = "High"
$ws.Range("J14").Value = "Yes"

# Row 15
$ws.Range("B15").Value = 11
$ws.Range("C15").Value = "Dashboard"
$ws.Range("D15").Value = "Functional Error"
$ws.Range("E15").Value = "Open"
$ws.Range("F15").Value = "AR"
$ws.Range("G15").Value = "Search bar does not return results"
$ws.Range("H15").Value = "Sacha Elkaim"
$ws.Range("I15").Value = "Medium"
$ws.Range("J15").Value = "No"

# Row 16
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = "Production - Inventory "
$ws.Range("D16").Value = "Functional Error"
$ws.Range("E16").Value = "Minor"
$ws.Range("F16").Value = "AR"
$ws.Range("G16").Value = "Search bar always returns Seat Rod, no matter search parameters"
$ws.Range("H16").Value = "Derek Ruiz-Cigana"
$ws.Range("I16").Value = "Medium"
$ws.Range("J16").Value = "No"

# Row 17
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "Admin Panel"
$ws.Range("D17").Value = "Functional Error"
$ws.Range("E17").Value = "Missing"
$ws.Range("F17").Value = "AR"
$ws.Range("G17").Value = "Cannot change a user's location when the option is selected"
$ws.Range("H17").Value = "Derek Ruiz-Cigana"
$ws.Range("I17").Value = "Medium"
$ws.Range("J17").Value = "No"

# Row 18
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = "Production - Inventory "
$ws.Range("D18").Value = "Functional Error"
$ws.Range("E18").Value = "Major"
$ws.Range("F18").Value = "AR"
$ws.Range("G18").Value = "Error using the create saddle button -- no end time"
$ws.Range("H18").Value = "Derek Ruiz-Cigana"
$ws.Range("I18").Value = "High"
$ws.Range("J18").Value = "No"
